$wb = $excel.ActiveWorkbook

# --- ALC row 12 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(12, 8).Value = 3897.6924
$ws.Cells.Item(12, 10).Value = 4500.1665
$ws.Cells.Item(12, 12).Value = 4500.1665
$ws.Cells.Item(12, 14).Value = -4840.1665

# --- ALC row 19 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value = 2567.2942
$ws.Cells.Item(19, 9).Value = 2107.6667
$ws.Cells.Item(19, 10).Value = 2818
$ws.Cells.Item(19, 11).Value = 2107.6667
$ws.Cells.Item(19, 12).Value = 2818
$ws.Cells.Item(19, 13).Value = -1932.6667
$ws.Cells.Item(19, 14).Value = -3168

# --- ALC row 33 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 386.65625
$ws.Cells.Item(33, 9).Value = 497.54544
$ws.Cells.Item(33, 11).Value = 497.54544
$ws.Cells.Item(33, 13).Value = -268.54544

# --- ALC row 98 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(98, 8).Value = 2362.2974
$ws.Cells.Item(98, 9).Value = 640.4783
$ws.Cells.Item(98, 10).Value = 5191
$ws.Cells.Item(98, 11).Value = 640.4783
$ws.Cells.Item(98, 12).Value = 5191
$ws.Cells.Item(98, 13).Value = 857.5217
$ws.Cells.Item(98, 14).Value = -8187

# --- ALC row 103 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(103, 8).Value = 425
$ws.Cells.Item(103, 10).Value = 460
$ws.Cells.Item(103, 12).Value = 1380
$ws.Cells.Item(103, 14).Value = -2552

# --- ALC row 122 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(122, 8).Value = 2362.2974
$ws.Cells.Item(122, 9).Value = 640.4783
$ws.Cells.Item(122, 10).Value = 5191
$ws.Cells.Item(122, 11).Value = 1921.4349
$ws.Cells.Item(122, 12).Value = 15573
$ws.Cells.Item(122, 13).Value = 528.5651
$ws.Cells.Item(122, 14).Value = -20473

# --- ALC row 127 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(127, 8).Value = 1050
$ws.Cells.Item(127, 9).Value = 1050
$ws.Cells.Item(127, 11).Value = 3150
$ws.Cells.Item(127, 13).Value = 1810

# --- ALC row 132 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(132, 8).Value = 4803.6763
$ws.Cells.Item(132, 9).Value = 4966.731
$ws.Cells.Item(132, 11).Value = 14900.193
$ws.Cells.Item(132, 13).Value = -12370.193

# --- ALC row 141 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(141, 8).Value = 4789.625
$ws.Cells.Item(141, 9).Value = 4789.625
$ws.Cells.Item(141, 11).Value = 14368.875
$ws.Cells.Item(141, 13).Value = -9188.875

# --- ARM row 32 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 1553.2892
$ws.Cells.Item(32, 9).Value = 1581.9012
$ws.Cells.Item(32, 11).Value = 1581.9012
$ws.Cells.Item(32, 13).Value = -1294.9012

# --- ARM row 74 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 12860
$ws.Cells.Item(74, 9).Value = 14713.238
$ws.Cells.Item(74, 11).Value = 14713.238
$ws.Cells.Item(74, 13).Value = -13839.238

# --- ARM row 77 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(77, 8).Value = 12860
$ws.Cells.Item(77, 9).Value = 14713.238
$ws.Cells.Item(77, 11).Value = 73566.19
$ws.Cells.Item(77, 13).Value = -69198.19

# --- ARM row 102 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(102, 8).Value = 3417.8333
$ws.Cells.Item(102, 9).Value = 2836.6667
$ws.Cells.Item(102, 10).Value = 3999
$ws.Cells.Item(102, 11).Value = 2836.6667
$ws.Cells.Item(102, 12).Value = 3999
$ws.Cells.Item(102, 13).Value = -1214.6667
$ws.Cells.Item(102, 14).Value = -7243

# --- ARM row 110 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(110, 8).Value = 115780.16
$ws.Cells.Item(110, 9).Value = 123838.17
$ws.Cells.Item(110, 11).Value = 123838.17
$ws.Cells.Item(110, 13).Value = -121793.17

# --- ARM row 122 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(122, 8).Value = 4732.3687
$ws.Cells.Item(122, 9).Value = 3800.1
$ws.Cells.Item(122, 11).Value = 11400.3
$ws.Cells.Item(122, 13).Value = -8950.299999999999

# --- ARM row 132 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 8060.6924
$ws.Cells.Item(132, 9).Value = 2787.25
$ws.Cells.Item(132, 11).Value = 8361.75
$ws.Cells.Item(132, 13).Value = -5831.75

# --- BSM row 20 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 1517.4117
$ws.Cells.Item(20, 9).Value = 1613.1111
$ws.Cells.Item(20, 10).Value = 1409.75
$ws.Cells.Item(20, 11).Value = 1613.1111
$ws.Cells.Item(20, 12).Value = 1409.75
$ws.Cells.Item(20, 13).Value = -1366.1111
$ws.Cells.Item(20, 14).Value = -1903.75

# --- BSM row 42 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(42, 8).Value = 179500
$ws.Cells.Item(42, 10).Value = 179500
$ws.Cells.Item(42, 12).Value = 179500
$ws.Cells.Item(42, 14).Value = -180156

# --- BSM row 43 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(43, 8).Value = 0
$ws.Cells.Item(43, 10).Value = 0
$ws.Cells.Item(43, 12).Value = 0
$ws.Cells.Item(43, 14).ClearContents()  # remove cell entirely (was -250046)

# --- BSM row 94 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 943.96295
$ws.Cells.Item(94, 9).Value = 668.1818
$ws.Cells.Item(94, 10).Value = 2157.4
$ws.Cells.Item(94, 11).Value = 668.1818
$ws.Cells.Item(94, 12).Value = 2157.4
$ws.Cells.Item(94, 13).Value = -217.1818
$ws.Cells.Item(94, 14).Value = -3059.4

# --- BSM row 99 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 1996
$ws.Cells.Item(99, 9).Value = 1374.5333
$ws.Cells.Item(99, 11).Value = 1374.5333
$ws.Cells.Item(99, 13).Value = 123.4666999999999

# --- BSM row 107 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 1505.2413
$ws.Cells.Item(107, 9).Value = 1389.091
$ws.Cells.Item(107, 10).Value = 1870.2858
$ws.Cells.Item(107, 11).Value = 1389.091
$ws.Cells.Item(107, 12).Value = 1870.2858
$ws.Cells.Item(107, 13).Value = 530.9090000000001
$ws.Cells.Item(107, 14).Value = -5710.2858

# --- BSM row 134 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 57009.74
$ws.Cells.Item(134, 9).Value = 2923.5
$ws.Cells.Item(134, 11).Value = 8770.5
$ws.Cells.Item(134, 13).Value = -6235.5

# --- BSM row 139 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(139, 8).Value = 99179.75
$ws.Cells.Item(139, 10).Value = 99179.75
$ws.Cells.Item(139, 12).Value = 99179.75
$ws.Cells.Item(139, 14).Value = -109459.75

# --- CRP row 16 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 4000
$ws.Cells.Item(16, 9).Value = 3500
$ws.Cells.Item(16, 11).Value = 3500
$ws.Cells.Item(16, 13).Value = -3213

# --- CRP row 31 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3173.65
$ws.Cells.Item(31, 9).Value = 1159.7778
$ws.Cells.Item(31, 11).Value = 1159.7778
$ws.Cells.Item(31, 13).Value = -864.7778000000001

# --- CRP row 34 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 3173.65
$ws.Cells.Item(34, 9).Value = 1159.7778
$ws.Cells.Item(34, 11).Value = 1159.7778
$ws.Cells.Item(34, 13).Value = -957.7778000000001

# --- CRP row 113 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(113, 8).Value = 4000
$ws.Cells.Item(113, 9).Value = 3500
$ws.Cells.Item(113, 11).Value = 3500
$ws.Cells.Item(113, 13).Value = -1330

# --- CRP row 122 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(122, 8).Value = 3493.6667
$ws.Cells.Item(122, 9).Value = 1756.3334
$ws.Cells.Item(122, 11).Value = 5269.0002
$ws.Cells.Item(122, 13).Value = -2819.0002

# --- CRP row 132 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 4074.375
$ws.Cells.Item(132, 9).Value = 3630.3845
$ws.Cells.Item(132, 11).Value = 10891.1535
$ws.Cells.Item(132, 13).Value = -8361.1535

# --- CRP row 134 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134, 8).Value = 671413
$ws.Cells.Item(134, 10).Value = 1254974.9
$ws.Cells.Item(134, 12).Value = 3764924.7
$ws.Cells.Item(134, 14).Value = -3769994.7

# --- CUL row 55 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(55, 8).Value = 38451.59
$ws.Cells.Item(55, 9).Value = 919.25
$ws.Cells.Item(55, 10).Value = 50000
$ws.Cells.Item(55, 11).Value = 2757.75
$ws.Cells.Item(55, 12).Value = 150000
$ws.Cells.Item(55, 13).Value = -2580.75
$ws.Cells.Item(55, 14).Value = -150354

# --- CUL row 97 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(97, 8).Value = 484.27274
$ws.Cells.Item(97, 9).Value = 522.6
$ws.Cells.Item(97, 10).Value = 452.33334
$ws.Cells.Item(97, 11).Value = 1567.8
$ws.Cells.Item(97, 12).Value = 1357.00002
$ws.Cells.Item(97, 13).Value = -1071.8
$ws.Cells.Item(97, 14).Value = -2349.00002

# --- CUL row 131 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 3598.7307
$ws.Cells.Item(131, 10).Value = 5223.3335
$ws.Cells.Item(131, 12).Value = 15670.0005
$ws.Cells.Item(131, 14).Value = -25750.0005

# --- GSM row 122 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 921570.5
$ws.Cells.Item(122, 9).Value = 1228115.5
$ws.Cells.Item(122, 10).Value = 1935.3334
$ws.Cells.Item(122, 11).Value = 3684346.5
$ws.Cells.Item(122, 12).Value = 5806.0002
$ws.Cells.Item(122, 13).Value = -3681896.5
$ws.Cells.Item(122, 14).Value = -10706.0002

# --- GSM row 126 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(126, 8).Value = 3631.5715
$ws.Cells.Item(126, 9).Value = 2128.2856
$ws.Cells.Item(126, 11).Value = 6384.8568
$ws.Cells.Item(126, 13).Value = -3914.8568

# --- LTW row 16 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 1044.2
$ws.Cells.Item(16, 9).Value = 938.55554
$ws.Cells.Item(16, 11).Value = 938.55554
$ws.Cells.Item(16, 13).Value = -768.55554

# --- LTW row 100 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(100, 8).Value = 152496.72
$ws.Cells.Item(100, 9).Value = 211380.4
$ws.Cells.Item(100, 10).Value = 5287.5
$ws.Cells.Item(100, 11).Value = 211380.4
$ws.Cells.Item(100, 12).Value = 5287.5
$ws.Cells.Item(100, 13).Value = -210839.4
$ws.Cells.Item(100, 14).Value = -6369.5

# --- LTW row 132 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 5405.8887
$ws.Cells.Item(132, 9).Value = 3034.6667
$ws.Cells.Item(132, 10).Value = 7777.1113
$ws.Cells.Item(132, 11).Value = 9104.000100000001
$ws.Cells.Item(132, 12).Value = 23331.3339
$ws.Cells.Item(132, 13).Value = -6574.000100000001
$ws.Cells.Item(132, 14).Value = -28391.3339

# --- WVR row 47 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(47, 8).Value = 13700
$ws.Cells.Item(47, 10).Value = 13700
$ws.Cells.Item(47, 12).Value = 13700
$ws.Cells.Item(47, 14).Value = -14844

# --- WVR row 104 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(104, 8).Value = 75000
$ws.Cells.Item(104, 10).Value = 75000
$ws.Cells.Item(104, 12).Value = 75000
$ws.Cells.Item(104, 14).Value = -81988

# --- WVR row 122 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 35718464
$ws.Cells.Item(122, 9).Value = 55558376
$ws.Cells.Item(122, 10).Value = 6617.8
$ws.Cells.Item(122, 11).Value = 166675128
$ws.Cells.Item(122, 12).Value = 19853.4
$ws.Cells.Item(122, 13).Value = -166672678
$ws.Cells.Item(122, 14).Value = -24753.4
